$wb = $excel.ActiveWorkbook

# Map of row -> new F value (想去人数 / "want to go" count)
$fUpdates = @{
    2  = 1122
    3  = 835
    5  = 50
    6  = 1110
    8  = 2058
    9  = 7677
    11 = 433
    12 = 362
    13 = 145
    14 = 410
    16 = 7846
    17 = 317
    18 = 1366
    22 = 162
    24 = 152
    26 = 23
    28 = 24
    29 = 419
    30 = 1140
    31 = 56
    33 = 65
    35 = 43
    37 = 69
}

# Sheets that contain this duplicated data table: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    # Row 17's lowest price (column G) changed from a numeric 45 to the
    # text "已售罄" (sold out)
    $ws.Cells.Item(17, 7).Value = "已售罄"
}
